$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two "FIXED $ OR $/UNIT/PERIOD" cells whose values changed:
#   D2: "$/period" -> "$/ft/period"
#   D4: "$/period" -> "$/sq.ft/period"
$ws.Range("D2").Value = "$/ft/period"
$ws.Range("D4").Value = "$/sq.ft/period"

# The whole C1:D15 block had an explicit (no-op) style applied (fontId 0,
# applyFont="1") that gets cleared back to the sheet's "Normal" style,
# letting the cells fall back to the column-level formatting again.
$ws.Range("C1:D15").Style = "Normal"

# Selection moved to D17 (single cell) before the file was saved.
$ws.Range("D17").Select()
